$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new exam entry for 2025 - Vår as row 18
$ws.Range("A18").Value = "2025 - Vår"
$ws.Range("B18").Value = "[Oppgaveformulering](tidligere-eksamensoppgaver/skole-25-v.pdf)"
$ws.Range("C18").Value = "[Løsningsforslag](tidligere-eksamensoppgaver/skole-25-v-fasit.pdf)"

# Update the selection to reflect the new active cell after the edit
$ws.Range("C19").Select()
